# Append new rows (166-178) of results data to the active worksheet,
# matching the "improve gui for buildings to" data update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(84.2,  0.36,                118.7564432621002,   111.640625),
    @(84.2,  0.36,                180.320063829422,    178.1875),
    @(39.78, 0.16,                32.275066614151,     30.96875),
    @(39.78, 0.16,                67.87305283546448,   49.28125),
    @(39.78, 0.16,                29.76566123962402,   28.8125),
    @(39.78, 0.16,                29.20565843582153,   28.984375),
    @(39.78, 0.16,                19.59266972541809,   19.453125),
    @(39.78, 0.16,                39.31414198875427,   39.421875),
    @(39.06, 0.04000000000000001, 5.450976371765137,   5.34375),
    @(39.06, 0.04000000000000001, 7.199516773223877,   7.1875),
    @(39.06, 0.04000000000000001, 6.611833810806274,   6.375),
    @(39.06, 0.04000000000000001, 5.285555362701416,   5.21875),
    @(39.78, 0.16,                19.69591951370239,   19.578125)
)

$startRow = 166
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
